$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "dropdown"
$ws.Range("F5").Value = "select"
$ws.Range("G5").Value = "None"
$ws.Range("H5").Value = "id,name,class etc"

$ws.Range("E6").Value = "item in dropdown"
$ws.Range("F6").Value = "option"
$ws.Range("G6").Value = "None"
$ws.Range("H6").Value = "id,name,class etc"

$ws.Range("H6").Select()
